# Revert "adding term 2.0.0"
# - Metadata sheet: roll Version/Date/Contact/Description back to the
#   pre-2.0.0 values (and restore the "interes" typo that 2.0.0 had fixed).
# - "Include from FSIII" sheet: drop the concept row that 2.0.0 added
#   (d7ff926a-4955-478f-b300-0b0ec0785013), shifting the rows below it up.
# - "Include from FFB" sheet is untouched.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "1.1.0"
$meta.Range("B8").Value = "2023-07-10T23:08:03+02:00"
$meta.Range("B10").Value = "No display for ContactDetail"
$meta.Range("B11").Value = "Matter of interes values to support when no observations have been made"

# --- Include from FSIII sheet ---------------------------------------
# Row 2 holds the concept added by 2.0.0; delete it so the following
# rows (B6 / blank / System URI row) shift up.
$fsiii = $wb.Worksheets.Item("Include from FSIII")
$fsiii.Rows.Item(2).Delete()
